$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 407 from 45175 to 45177
$ws.Range("C2:C407").Value = 45177
